$wb = $excel.ActiveWorkbook

# --- Sheet "info": update the recorded run time in B2 ---
$wsInfo = $wb.Worksheets.Item("info")
$wsInfo.Range("B2").Value = 15.58531665802002

# --- Sheet "x": update computed values in B2:B151 ---
$wsX = $wb.Worksheets.Item("x")
$xValues = @(
    0.8763319097879275, 0.04373774937504579, 0.03468373649165744, 0.04266841321593796, 0.03935354646499981, 1.009140262103619, 0.04118348011900557, 0.03675866401283336, 0.04104874882348094, 0.0314335017432567,
    0.02971150682067697, 0.03388049353486607, 0.0307769291871152, 0.03205347195563703, 0.02952935556634553, 0.0332600788682887, 0.03161328439572575, 0.03158319606982835, 0.03024760936101478, 0.02982463242837517,
    0.03073488009708395, 1.019712998852669, 0.03953734820436428, 0.03934137311963928, 0.03847993546200173, 0.03025974111604739, 0.0353032813806186, 1.016175372726841, 0.02985238699994165, 0.03601819190285348,
    0.03864539704629352, 0.03800628460852666, 0.03509723663459581, 1.018929841724191, 0.03625975147309615, 0.02934917141790918, 0.03551075787861446, 0.03967085542826821, 0.02992185128470725, 0.03424390906740583,
    0.0306817474428441, 0.03770938339796523, 1.010095320970243, 0.02972347852863355, 0.03866187475220857, 0.03897382643452863, 0.02838855691571966, 0.02961097392394872, 0.03104262803324215, 0.03579070431973548,
    0.03268588703078296, 0.02992732548652394, 0.03177746167622392, 0.03784212831429859, 0.03022270705672732, 0.03019029651266226, 0.03209441900059744, 0.03086680935632359, 0.02812851790354326, 0.02832829234238862,
    0.0320647828112539, 0.03576843225669901, 0.02983829538716302, 0.04035443593752901, 0.03001774524176469, 0.02918166623816583, 0.03454149924869043, 0.04400849107865243, 0.02827952430562559, 0.03488983752281195,
    0.03384071679736969, 0.02825015273484522, 0.03228728367407334, 1.019074349557362, 0.02915962735046505, 0.04415846929723287, 0.03072508887920293, 0.02900113224157901, 0.03064900109360261, 0.03516727143937166,
    0.03938955251311569, 0.02913842049286626, 0.03053731574781805, 0.03933378526954202, 0.03033429025573767, 0.0396293645455421, 0.03705090657118943, 0.03097799132157467, 0.03555782499631806, 0.0356529618729443,
    0.03489288481094, 0.03906102735963154, 0.03529817091611722, 0.03894114387198945, 0.02869592438405543, 0.03507889105881718, 0.03065002040207827, 0.03617876503603727, 0.03996763219164487, 0.0405211675160703,
    0.03081781060878424, 0.030833497154624, 0.02954309312721515, 0.0297073704817199, 0.03928764563552023, 0.03205414094800086, 0.03225418288526143, 0.0295651963481484, 0.03016783028002062, 0.03284508073523334,
    0.03419195001846054, 0.03367906727983902, 0.03376968220804168, 0.02873304872721034, 0.03200279060242, 0.03097969423667087, 0.03493235250776237, 0.04287782246394797, 1.012547060783924, 0.03073130303091414,
    0.03126679214639315, 0.03002069488185221, 0.03927832600249242, 0.03562826848532147, 0.03191824786638869, 0.03688389201551465, 0.03778625086173282, 1.019593833493061, 0.03030990053322799, 0.02814327915421121,
    0.03704024986978131, 0.03283815618271361, 0.03655472480055769, 0.03946629169899457, 0.03438809215665058, 0.03890110995547402, 0.02928562797624953, 0.04289718703050781, 0.03251592761243068, 0.02916325765567453,
    0.04011185500743412, -0.9036523360630978, 0.03031549225566159, 0.02886302873755588, 0.03625979683875341, 0.0340165205695327, 0.02993868421867408, 0.02935892987029709, 0.03606804604411745, 0.02836311206360224
)

$xArr = New-Object 'object[,]' $xValues.Length,1
for ($i = 0; $i -lt $xValues.Length; $i++) {
    $xArr[$i,0] = $xValues[$i]
}
$wsX.Range("B2:B151").Value = $xArr

# --- Sheet "y": reset all indicator values in B2:B151 to 0 ---
$wsY = $wb.Worksheets.Item("y")
$yArr = New-Object 'object[,]' 150,1
for ($i = 0; $i -lt 150; $i++) {
    $yArr[$i,0] = 0
}
$wsY.Range("B2:B151").Value = $yArr
